# Auto-generated edit script: update F-column ('想去人数' / want-to-go count) values
# across all four worksheets to match the refreshed gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 139
$ws.Cells.Item(4, 6).Value = 798
$ws.Cells.Item(6, 6).Value = 633
$ws.Cells.Item(7, 6).Value = 1196
$ws.Cells.Item(8, 6).Value = 98
$ws.Cells.Item(9, 6).Value = 769
$ws.Cells.Item(10, 6).Value = 677
$ws.Cells.Item(11, 6).Value = 254
$ws.Cells.Item(13, 6).Value = 351
$ws.Cells.Item(15, 6).Value = 872
$ws.Cells.Item(16, 6).Value = 9610
$ws.Cells.Item(17, 6).Value = 581
$ws.Cells.Item(18, 6).Value = 581
$ws.Cells.Item(21, 6).Value = 40
$ws.Cells.Item(23, 6).Value = 253
$ws.Cells.Item(24, 6).Value = 1731
$ws.Cells.Item(26, 6).Value = 281
$ws.Cells.Item(30, 6).Value = 254
$ws.Cells.Item(31, 6).Value = 181
$ws.Cells.Item(32, 6).Value = 257
$ws.Cells.Item(33, 6).Value = 54
$ws.Cells.Item(37, 6).Value = 176
$ws.Cells.Item(38, 6).Value = 157

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 33
$ws.Cells.Item(6, 6).Value = 104
$ws.Cells.Item(7, 6).Value = 113
$ws.Cells.Item(10, 6).Value = 229
$ws.Cells.Item(11, 6).Value = 4432
$ws.Cells.Item(12, 6).Value = 74
$ws.Cells.Item(15, 6).Value = 58
$ws.Cells.Item(16, 6).Value = 257

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 802

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 802
$ws.Cells.Item(5, 6).Value = 139
$ws.Cells.Item(6, 6).Value = 798
$ws.Cells.Item(7, 6).Value = 33
$ws.Cells.Item(9, 6).Value = 633
$ws.Cells.Item(10, 6).Value = 1196
$ws.Cells.Item(12, 6).Value = 104
$ws.Cells.Item(13, 6).Value = 113
$ws.Cells.Item(14, 6).Value = 769
$ws.Cells.Item(15, 6).Value = 677
$ws.Cells.Item(16, 6).Value = 254
$ws.Cells.Item(17, 6).Value = 351
$ws.Cells.Item(19, 6).Value = 872
$ws.Cells.Item(20, 6).Value = 9610
$ws.Cells.Item(21, 6).Value = 229
$ws.Cells.Item(22, 6).Value = 582
$ws.Cells.Item(23, 6).Value = 582
$ws.Cells.Item(25, 6).Value = 40
$ws.Cells.Item(26, 6).Value = 253
$ws.Cells.Item(27, 6).Value = 1731
$ws.Cells.Item(29, 6).Value = 281
$ws.Cells.Item(32, 6).Value = 74
$ws.Cells.Item(33, 6).Value = 74
$ws.Cells.Item(36, 6).Value = 58
$ws.Cells.Item(37, 6).Value = 254
$ws.Cells.Item(38, 6).Value = 181
$ws.Cells.Item(39, 6).Value = 257
$ws.Cells.Item(40, 6).Value = 54
$ws.Cells.Item(47, 6).Value = 176
$ws.Cells.Item(48, 6).Value = 157
